$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 22500.1
$ws.Range("I51").Value = 20624.875
$ws.Range("J51").Value = 30001
$ws.Range("K51").Value = 20624.875
$ws.Range("L51").Value = 30001
$ws.Range("M51").Value = -20140.875
$ws.Range("N51").Value = -30969

# Row 61
$ws.Range("H61").Value = 2648.4285
$ws.Range("I61").Value = 408.8
$ws.Range("J61").Value = 8247.5
$ws.Range("K61").Value = 1226.4
$ws.Range("L61").Value = 24742.5
$ws.Range("M61").Value = -1054.4
$ws.Range("N61").Value = -25086.5

# Row 82
$ws.Range("H82").Value = 3031.75
$ws.Range("I82").Value = 1398.3636
$ws.Range("K82").Value = 4195.0908
$ws.Range("M82").Value = -3789.0908

# Row 85
$ws.Range("H85").Value = 3031.75
$ws.Range("I85").Value = 1398.3636
$ws.Range("K85").Value = 4195.0908
$ws.Range("M85").Value = -2791.0908

# Row 112
$ws.Range("H112").Value = 2447.7334
$ws.Range("I112").Value = 1474.5
$ws.Range("J112").Value = 2597.4614
$ws.Range("K112").Value = 4423.5
$ws.Range("L112").Value = 7792.3842
$ws.Range("M112").Value = -3315.5
$ws.Range("N112").Value = -10008.3842

# Row 113
$ws.Range("H113").Value = 6669.357
$ws.Range("J113").Value = 6276
$ws.Range("L113").Value = 6276
$ws.Range("N113").Value = -12784

# Row 135
$ws.Range("H135").Value = 690
$ws.Range("I135").Value = 713
$ws.Range("J135").Value = 506
$ws.Range("K135").Value = 6417
$ws.Range("L135").Value = 4554
$ws.Range("M135").Value = -3882
$ws.Range("N135").Value = -9624

# Row 137
$ws.Range("H137").Value = 1084.5
$ws.Range("I137").Value = 997.6429000000001
$ws.Range("J137").Value = 2300.5
$ws.Range("K137").Value = 2992.9287
$ws.Range("L137").Value = 6901.5
$ws.Range("M137").Value = -442.9287000000004
$ws.Range("N137").Value = -12001.5

# Row 138
$ws.Range("H138").Value = 3781.0264
$ws.Range("J138").Value = 3741.037
$ws.Range("L138").Value = 11223.111
$ws.Range("N138").Value = -21503.111

$ws = $wb.Worksheets.Item("ARM")
# Row 103
$ws.Range("H103").Value = 100000
$ws.Range("J103").Value = 100000
$ws.Range("L103").Value = 100000
$ws.Range("N103").Value = -102344

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 2799.6
$ws.Range("I64").Value = 2298.75
$ws.Range("J64").Value = 3372
$ws.Range("K64").Value = 2298.75
$ws.Range("L64").Value = 3372
$ws.Range("M64").Value = -2073.75
$ws.Range("N64").Value = -3822

# Row 67
$ws.Range("H67").Value = 2799.6
$ws.Range("I67").Value = 2298.75
$ws.Range("J67").Value = 3372
$ws.Range("K67").Value = 2298.75
$ws.Range("L67").Value = 3372
$ws.Range("M67").Value = -1518.75
$ws.Range("N67").Value = -4932

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 107
$ws.Range("H107").Value = 1493.1111
$ws.Range("I107").Value = 1204.75
$ws.Range("K107").Value = 1204.75
$ws.Range("M107").Value = 715.25

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 133
$ws.Range("H133").Value = 99998.664
$ws.Range("J133").Value = 99998.664
$ws.Range("L133").Value = 99998.664
$ws.Range("N133").Value = -110118.664

# Row 134
$ws.Range("H134").Value = 2520.7
$ws.Range("I134").Value = 2399.2856
$ws.Range("K134").Value = 7197.8568
$ws.Range("M134").Value = -4662.8568

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2388.5
$ws.Range("I31").Value = 1790.25
$ws.Range("K31").Value = 1790.25
$ws.Range("M31").Value = -1495.25

# Row 34
$ws.Range("H34").Value = 2388.5
$ws.Range("I34").Value = 1790.25
$ws.Range("K34").Value = 1790.25
$ws.Range("M34").Value = -1588.25

# Row 99
$ws.Range("H99").Value = 1703.8182
$ws.Range("I99").Value = 1363.8572
$ws.Range("J99").Value = 2298.75
$ws.Range("K99").Value = 1363.8572
$ws.Range("L99").Value = 2298.75
$ws.Range("M99").Value = 134.1428000000001
$ws.Range("N99").Value = -5294.75

# Row 107
$ws.Range("H107").Value = 2459.2917
$ws.Range("I107").Value = 564.1429000000001
$ws.Range("K107").Value = 564.1429000000001
$ws.Range("M107").Value = 1355.8571

# Row 126
$ws.Range("H126").Value = 1703.8182
$ws.Range("I126").Value = 1363.8572
$ws.Range("J126").Value = 2298.75
$ws.Range("K126").Value = 4091.5716
$ws.Range("L126").Value = 6896.25
$ws.Range("M126").Value = -1621.5716
$ws.Range("N126").Value = -11836.25

# Row 134
$ws.Range("H134").Value = 26551.293
$ws.Range("I134").Value = 27759.076
$ws.Range("K134").Value = 83277.228
$ws.Range("M134").Value = -80742.228

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 1696
$ws.Range("I3").Value = 1696
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5088
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4976
$ws.Range("N3").ClearContents()

# Row 68
$ws.Range("H68").Value = 4198
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4198
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 12594
$ws.Range("N68").Value = -14216
$ws.Range("M68").ClearContents()

# Row 71
$ws.Range("H71").Value = 4198
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4198
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 37782
$ws.Range("N71").Value = -45894
$ws.Range("M71").ClearContents()

# Row 131
$ws.Range("H131").Value = 20574.875
$ws.Range("J131").Value = 20467.637
$ws.Range("L131").Value = 61402.91099999999
$ws.Range("N131").Value = -71482.91099999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# Row 46
$ws.Range("H46").Value = 2371.1667
$ws.Range("I46").Value = 1349.5
$ws.Range("K46").Value = 1349.5
$ws.Range("M46").Value = -1161.5

# Row 122
$ws.Range("H122").Value = 3222.9
$ws.Range("I122").Value = 3077.1428
$ws.Range("K122").Value = 9231.428400000001
$ws.Range("M122").Value = -6781.428400000001

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 1512000
$ws.Range("I29").Value = 3000000
$ws.Range("K29").Value = 3000000
$ws.Range("M29").Value = -2999710
